# Update gh-pages to output generated at 456a3b4
# Refreshes the "展览" (Exhibition) and "全部类型" (All types) sheets with
# newly scraped listing data:
#   - "展览": row 22 was a stray duplicate of row 21 and is removed outright,
#     shifting rows 23-43 up one slot (row 43 disappears) and picking up
#     several new listings along the way; several "想去人数" (want-to-go
#     count) figures were also refreshed.
#   - "全部类型": same "想去人数" figures refreshed (no row shift there).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览": delete the duplicate row 22 (Aw7th duplicate of row 21),
# which shifts all subsequent rows up by one and drops the old row 43.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Rows.Item(22).Delete() | Out-Null

# Column A holds a plain 0-based sequence number (row - 1); restore it for
# every data row now that the rows shifted up.
$lastRow = $wsExpo.Cells.Item($wsExpo.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $wsExpo.Cells.Item($r, 1).Value = $r - 2
}

# Refreshed "想去人数" (F column) counts at their final (post-delete) rows.
$expoFUpdates = @{
    4  = 3768
    5  = 298
    6  = 5315
    7  = 600
    8  = 435
    10 = 1069
    14 = 735
    16 = 49
    21 = 6094
    25 = 7038
    29 = 375
    30 = 757
    35 = 1172
    40 = 1159
}
foreach ($row in $expoFUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoFUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "全部类型": same refreshed "想去人数" counts, no row shift here.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allFUpdates = @{
    7  = 3768
    8  = 298
    9  = 5315
    10 = 600
    11 = 435
    13 = 1069
    17 = 735
    19 = 49
    25 = 6094
    29 = 7038
    33 = 375
    34 = 757
    40 = 1172
    45 = 1159
}
foreach ($row in $allFUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allFUpdates[$row]
}
